$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Uppercase all the names in column A (rows 2-11)
$ws.Range("A2").Value = "ARTHUR"
$ws.Range("A3").Value = "BEATRIZ"
$ws.Range("A4").Value = "CARLOS"
$ws.Range("A5").Value = "DANIELA"
$ws.Range("A6").Value = "EDUARDO"
$ws.Range("A7").Value = "FERNANDA"
$ws.Range("A8").Value = "GUSTAVO"
$ws.Range("A9").Value = "HUGO"
$ws.Range("A10").Value = "ISABELA"
$ws.Range("A11").Value = "JOÃO"

# Update the selection/active cell to A2
$ws.Range("A2").Select()
